# Update database and change read_price algorithm
# - drops the oldest fiscal-period column's data, shifts remaining periods
#   one column to the left, and appends a new period (1401/12) in column H
# - updates the "publish date" row to match
# - row 15 (previously displayed as "-") is now a real 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers (D:H) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates (D:H) ---
$ws.Range("D9").Value = "1399-02-31 (10)"
$ws.Range("E9").Value = "1400-02-30 (8)"
$ws.Range("F9").Value = "1401-02-31 (8)"
$ws.Range("G9").Value = "1402-02-10 (7)"
# H9 ("1402-02-10") looks like an ISO date, so a plain .Value assignment
# would be auto-converted into a date serial number. Route it through a
# text formula and paste the computed value back so it lands as a plain
# text cell, matching the original cell's shape/styling.
$ws.Range("H9").Formula = '="1402-02-10"'
$ws.Range("H9").Copy()
$ws.Range("H9").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 9508850
$ws.Range("E11").Value = 14124456
$ws.Range("F11").Value = 22157318
$ws.Range("G11").Value = 52733742
$ws.Range("H11").Value = 79362180

# --- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) ---
$ws.Range("D12").Value = -3692535
$ws.Range("E12").Value = -5544019
$ws.Range("F12").Value = -7370112
$ws.Range("G12").Value = -18866925
$ws.Range("H12").Value = -24756075

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = 5816315
$ws.Range("E13").Value = 8580437
$ws.Range("F13").Value = 14787206
$ws.Range("G13").Value = 33866817
$ws.Range("H13").Value = 54606105

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ---
$ws.Range("D14").Value = -906563
$ws.Range("E14").Value = -811108
$ws.Range("F14").Value = -1212981
$ws.Range("G14").Value = -2088702
$ws.Range("H14").Value = -4955612

# --- Row 15: هزینه کاهش ارزش دریافتنی‌‏ها (was text "-", now numeric 0) ---
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---
$ws.Range("D16").Value = 355205
$ws.Range("E16").Value = 334587
$ws.Range("F16").Value = 168096
$ws.Range("G16").Value = 473189
$ws.Range("H16").Value = 1792367

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Range("D17").Value = 5264957
$ws.Range("E17").Value = 8103916
$ws.Range("F17").Value = 13742321
$ws.Range("G17").Value = 32251304
$ws.Range("H17").Value = 51442860

# --- Row 18: هزینه های مالی (Financial expenses) -- stays all zero ---
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ---
$ws.Range("D19").Value = 681357
$ws.Range("E19").Value = 1073090
$ws.Range("F19").Value = 1470153
$ws.Range("G19").Value = 2431133
$ws.Range("H19").Value = 4036071

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---
$ws.Range("D20").Value = 5946314
$ws.Range("E20").Value = 9177006
$ws.Range("F20").Value = 15212474
$ws.Range("G20").Value = 34682437
$ws.Range("H20").Value = 55478931

# --- Row 21: مالیات (Tax) ---
$ws.Range("D21").Value = -356834
$ws.Range("E21").Value = -426882
$ws.Range("F21").Value = -845462
$ws.Range("G21").Value = -2227521
$ws.Range("H21").Value = -6008416

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ---
$ws.Range("D22").Value = 5589480
$ws.Range("E22").Value = 8750124
$ws.Range("F22").Value = 14367012
$ws.Range("G22").Value = 32454916
$ws.Range("H22").Value = 49470515

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی -- stays all zero ---
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# --- Row 24: سود (زیان) خالص (Net profit) ---
$ws.Range("D24").Value = 5589480
$ws.Range("E24").Value = 8750124
$ws.Range("F24").Value = 14367012
$ws.Range("G24").Value = 32454916
$ws.Range("H24").Value = 49470515

# --- Row 25: سود هر سهم پس از کسر مالیات (EPS after tax) ---
$ws.Range("D25").Value = 3123
$ws.Range("E25").Value = 4889
$ws.Range("F25").Value = 8027
$ws.Range("G25").Value = 18132
$ws.Range("H25").Value = 27639

# --- Row 26: سرمایه (Capital) -- stays constant ---
$ws.Range("D26").Value = 1789912
$ws.Range("E26").Value = 1789912
$ws.Range("F26").Value = 1789912
$ws.Range("G26").Value = 1789912
$ws.Range("H26").Value = 1789912

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه ---
$ws.Range("D27").Value = 3123
$ws.Range("E27").Value = 4889
$ws.Range("F27").Value = 8027
$ws.Range("G27").Value = 18132
$ws.Range("H27").Value = 27639
